# Verbs workbook update:
#  - add the missing "Ba Form" (conditional) readings for several verbs
#    that were already in the sheet (rows 9-15: 切る/知る/飛ぶ/座る/乗る/撮る/聞く)
#  - add a brand-new verb row (92: 稼ぐ - "to earn") with all six other
#    conjugation forms plus its own Ba Form
#  - move the active selection to the newly important cell (G15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New row for 稼ぐ (kasegu, "to earn") -----------------------------
$ws.Range("A92:G92").Font.Name = "Yu Gothic"
$ws.Range("A92").Value = "稼ぐ"
$ws.Range("B92").Value = "稼いで"
$ws.Range("C92").Value = "稼いだ"
$ws.Range("D92").Value = "稼がない"
$ws.Range("E92").Value = "稼ぎます"
$ws.Range("F92").Value = "稼ごう"
$ws.Range("G92").Value = "稼げば"
$ws.Rows.Item(92).RowHeight = 18.75

# --- 2. Fill in the previously-blank "Ba Form" column for existing verbs ---
# Some of these cells also pick up the Yu Gothic font used elsewhere in
# that column; two (G13, G15) were already on the other CJK font and keep it.
$ws.Range("G9").Font.Name = "Yu Gothic"
$ws.Range("G9").Value = "切れば"

$ws.Range("G10").Font.Name = "Yu Gothic"
$ws.Range("G10").Value = "知れば"

$ws.Range("G11").Font.Name = "Yu Gothic"
$ws.Range("G11").Value = "飛べば"

$ws.Range("G12").Font.Name = "Yu Gothic"
$ws.Range("G12").Value = "座れば"

$ws.Range("G13").Value = "乗れば"

$ws.Range("G14").Font.Name = "Yu Gothic"
$ws.Range("G14").Value = "撮れば"

$ws.Range("G15").Value = "聞けば"

# --- 3. Update the saved selection/view -----------------------------------
$ws.Range("G15").Select()

Write-Output "verbs workbook updated"
